$d = $word.ActiveDocument

$replacements = @(
    @("21×84=", "53×63="),
    @("47×30=", "62×46="),
    @("45×49=", "89×96="),
    @("47×51=", "67×82="),
    @("97×56=", "44×29="),
    @("55×36=", "27×21="),
    @("40×65=", "35×55="),
    @("16×33=", "49×60="),
    @("54×83=", "97×76="),
    @("28×79=", "94×19="),
    @("14×60=", "53×33="),
    @("51×69=", "47×71="),
    @("37×25=", "75×87="),
    @("87×24=", "52×40="),
    @("28×38=", "47×94="),
    @("28×45=", "34×86="),
    @("29×32=", "57×26="),
    @("58×26=", "89×80="),
    @("33×34=", "81×79="),
    @("27×17=", "53×63="),
    @("88×95=", "81×11="),
    @("86×51=", "18×33="),
    @("64×71=", "36×92="),
    @("77×20=", "81×98="),
    @("83×94=", "36×80=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
